# LOM3037.xlsx content fix-up:
#  - "Objetivos:" row gets its real objectives text (was showing a teacher name)
#  - the 3 responsible-teacher names get their own rows (13-15) under
#    "Docentes responsáveis:" instead of being smeared across other rows
#  - every row below shifts down by 3 and picks up its real text
#    (short syllabus, full syllabus, evaluation method/criteria/recovery
#    rule, bibliography) instead of the misplaced teacher/criteria text
#    that used to sit there

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 3 new "teacher" rows right after "Docentes responsáveis:" (row 12)
$ws.Range("A13:A15").EntireRow.Insert()

# The insert drags a formatted-but-empty A column cell into the new rows;
# drop it so rows 13-15 end up with no <c> in column A, matching the target.
$ws.Range("A13:A15").Clear()

# Pick up the normal-column-B / red-column-C look from a row that already
# has it (row 16, "Programa resumido:") instead of re-deriving it, which
# avoids picking up column A's bold default for the still-empty B/C cells.
$ws.Range("B16").Copy()
$ws.Range("B13:B15").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C13:C15").PasteSpecial(-4122)

# ---- Objetivos: (row 10) now holds the real objectives paragraph ----
$ws.Range("B10").Value = "- Capacitar o aluno para relacionar as propriedades químicas e físicas dos elementos e seus compostos com suas posições na tabela periódica.-Capacitar o aluno a escrever os métodos industriais de obtenção dos elementos e seus compostos, bem como descrever suas aplicações- capacitar o aluno a comunicar-se eficazmente nas formas escrita, oral e gráfica"
$ws.Range("C10").Value = "- Capacitar o aluno para relacionar as propriedades químicas e físicas dos elementos e seus compostos com suas posições na tabela periódica.-Capacitar o aluno a escrever os métodos industriais de obtenção dos elementos e seus compostos, bem como descrever suas aplicações- capacitar o aluno a comunicar-se eficazmente nas formas escrita, oral e gráfica"

# ---- 3 new teacher rows (13-15) ----
$ws.Range("B13").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("C13").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("B14").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("C14").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"

# ---- Programa resumido: (row 16) ----
$ws.Range("B16").Value = "- Metais Representativos: Metais do Grupo 1, Metais do Grupo 2 e Metais do Grupo 13.- Metais de Transição: Propriedades gerais, Complexos."
$ws.Range("C16").Value = "- Metais Representativos: Metais do Grupo 1, Metais do Grupo 2 e Metais do Grupo 13.- Metais de Transição: Propriedades gerais, Complexos."

# ---- Programa: (row 18) ----
$ws.Range("B18").Value = "Metais e compostos dos grupos 1, 2, 13 e de transição da Tabela Periódica: Propriedades físicas e químicas (relação com a posição na Tabela Periódica), processos de obtenção dos metais e compostos e aplicações - Formação de Complexos.Relacionar a disciplina com disciplinas anteriores e posteriores da grade do curso."
$ws.Range("C18").Value = "Metais e compostos dos grupos 1, 2, 13 e de transição da Tabela Periódica: Propriedades físicas e químicas (relação com a posição na Tabela Periódica), processos de obtenção dos metais e compostos e aplicações - Formação de Complexos.Relacionar a disciplina com disciplinas anteriores e posteriores da grade do curso."

# ---- Método: (row 21) ----
$ws.Range("B21").Value = "A avaliação tem como requisito quantificar as competências adquiridas conforme objetivadas.Duas provas escritas (P1 e P2) e listas de exercícios de acompanhamento continuado. A partir das notas das listas de exercício será calculada a média, LE."
$ws.Range("C21").Value = "A avaliação tem como requisito quantificar as competências adquiridas conforme objetivadas.Duas provas escritas (P1 e P2) e listas de exercícios de acompanhamento continuado. A partir das notas das listas de exercício será calculada a média, LE."

# ---- Critério: (row 22) ----
$ws.Range("B22").Value = "NF = (P1 + P2 + LE) /3"
$ws.Range("C22").Value = "NF = (P1 + P2 + LE) /3"

# ---- Norma de recuperação: (row 23) ----
$ws.Range("B23").Value = "Será realizada uma prova escrita valendo de zero a dez (NR) e a média final calculada pela equação: NF + NR"
$ws.Range("C23").Value = "Será realizada uma prova escrita valendo de zero a dez (NR) e a média final calculada pela equação: NF + NR"

# ---- Bibliografia: (row 24) ----
$ws.Range("B24").Value = "LEE, J. D. “Química Inorgânica não tão Concisa”, Editora Edgard Blücher, 1999. - SHRIVER, D. F.; ATKINS, P. W. “Química Inorgânica”, Editora Bookman, 4ª edição, 2008. - QUAGLIANO, J. V.; VALLARINO, L. “Química”, Editora Guanabara Koogan, 1973. - BUCHEL, K. H.; MORETTO, H. H.; WODITSCH, P. “Industrial Inorganic Chemistry”, Editora Wiley-VCH, 2000. - RAYNER-CANHAM, G.; OVERTON, T. “Química Inorgânica Descritiva”, Editora: Gen-LTC, 5ª edição, 2015. - SOUZA, M.M.V.M. “Processos Inorgânicos”, Editora: Synergia, 1ª edição, 2012."
$ws.Range("C24").Value = "LEE, J. D. “Química Inorgânica não tão Concisa”, Editora Edgard Blücher, 1999. - SHRIVER, D. F.; ATKINS, P. W. “Química Inorgânica”, Editora Bookman, 4ª edição, 2008. - QUAGLIANO, J. V.; VALLARINO, L. “Química”, Editora Guanabara Koogan, 1973. - BUCHEL, K. H.; MORETTO, H. H.; WODITSCH, P. “Industrial Inorganic Chemistry”, Editora Wiley-VCH, 2000. - RAYNER-CANHAM, G.; OVERTON, T. “Química Inorgânica Descritiva”, Editora: Gen-LTC, 5ª edição, 2015. - SOUZA, M.M.V.M. “Processos Inorgânicos”, Editora: Synergia, 1ª edição, 2012."

$excel.CutCopyMode = $false

Write-Output "edit complete"
